$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from an existing header cell (e.g. E1) to the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag data for columns F (KNN), G (SVM), H (RF)
$values = @(
  @($false, $false, $false),  # row 2
  @($false, $false, $false),  # row 3
  @($false, $false, $false),  # row 4
  @($false, $false, $false),  # row 5
  @($false, $false, $false),  # row 6
  @($false, $false, $false),  # row 7
  @($false, $false, $false),  # row 8
  @($true,  $false, $true),   # row 9
  @($true,  $true,  $true),   # row 10
  @($false, $false, $false),  # row 11
  @($false, $false, $false),  # row 12
  @($false, $false, $false),  # row 13
  @($false, $false, $false),  # row 14
  @($false, $false, $false),  # row 15
  @($false, $false, $false),  # row 16
  @($false, $false, $false),  # row 17
  @($true,  $true,  $true),   # row 18
  @($false, $false, $false),  # row 19
  @($false, $false, $false),  # row 20
  @($false, $false, $true),   # row 21
  @($false, $false, $false),  # row 22
  @($false, $false, $false),  # row 23
  @($false, $false, $false),  # row 24
  @($true,  $true,  $true)    # row 25
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $r = $i + 2
  $row = $values[$i]
  $ws.Cells.Item($r, 6).Value = $row[0]
  $ws.Cells.Item($r, 7).Value = $row[1]
  $ws.Cells.Item($r, 8).Value = $row[2]
}
